$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column A (Artista) - new ranking / names for rows 2..21
# ---------------------------------------------------------------------------
$names = @(
  "Marília Mendonça",
  "Bruno e Marrone",
  "Henrique e Juliano",
  "Simone e Simaria",
  "Raça Negra",
  "Alok",
  "Zé Neto e Cristiano",
  "Jorge e Mateus",
  "Péricles",
  "Sorriso Maroto",
  "Thiaguinho",
  "Maiara e Maraisa",
  " Os Barões da Pisadinha",
  "Marcos e Belutti",
  "Mano Walter",
  "Wesley Safadão",
  "Ferrugem",
  "Xand Avião",
  "Gusttavo Lima",
  "César Menotti e Fabiano"
)

for ($i = 0; $i -lt $names.Length; $i++) {
  $row = $i + 2
  $ws.Range("A$row").Value = $names[$i]
}

# ---------------------------------------------------------------------------
# Column B (Visualização) - updated view counts for rows 2..21
# ---------------------------------------------------------------------------
$views = @{
  2  = 55000000
  3  = 30000000
  4  = 21000000
  5  = 18000000
  6  = 16000000
  7  = 15000000
  8  = 14000000
  9  = 13000000
  10 = 11000000
  11 = 10000000
  12 = 10000000
  13 = 9000000
  14 = 8300000
  15 = 8000000
  16 = 7800000
  17 = 6700000
  18 = 6700000
  19 = 6200000
  20 = 6200000
  21 = 6000000
}

foreach ($row in $views.Keys) {
  $ws.Range("B$row").Value = $views[$row]
}

# ---------------------------------------------------------------------------
# Cell-specific formatting changes
# ---------------------------------------------------------------------------

# B12 gains right horizontal alignment (keeps default font)
$b12 = $ws.Range("B12")
$b12.HorizontalAlignment = -4152

# A15 gains a left/center aligned, black "arial" font
$a15 = $ws.Range("A15")
$a15.Font.ColorIndex = 1
$a15.Font.Name = "arial"
$a15.HorizontalAlignment = -4131
$a15.VerticalAlignment = -4108

# B15 gains a right/center aligned, black "Arial" font
$b15 = $ws.Range("B15")
$b15.Font.ColorIndex = 1
$b15.HorizontalAlignment = -4152
$b15.VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# View state: zoom + active selection
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 140
$ws.Range("A24:B24").Select()

Write-Host "Lives data updated"
